# Update 16 May: add new uploads
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "રણજીતભાઇ વાળા"
$ws.Range("B3").Value = "હિરેન સોજીત્રા"
